$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Mark1"
$ws.Range("C1").Value = "Mark2"

$ws.Range("A2").Value = "Indhu"
$ws.Range("B2").Value = 58
$ws.Range("C2").Value = 45

$ws.Range("A3").Value = "Lucky"
$ws.Range("B3").Value = 58
$ws.Range("C3").Value = 74

$ws.Range("A4").Value = "Leo"
$ws.Range("B4").Value = 58
$ws.Range("C4").Value = 34

$ws.Range("A5").Value = "Luna"
$ws.Range("B5").Value = 58
$ws.Range("C5").Value = 55

$ws.Range("C6").Select()
